$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "62.764.90"
Set-TextValue $ws.Range("E2") "  -1.61%  "
Set-TextValue $ws.Range("D3") "3.230.37"
Set-TextValue $ws.Range("E3") "  -1.78%  "
Set-TextValue $ws.Range("E4") "  +0.05%  "
Set-TextValue $ws.Range("D5") "594.19"
Set-TextValue $ws.Range("E5") "  -0.85%  "
Set-TextValue $ws.Range("D6") "136.47"
Set-TextValue $ws.Range("E6") "  -4.92%  "
Set-TextValue $ws.Range("E7") "  +0.11%  "
Set-TextValue $ws.Range("D8") "3.230.63"
Set-TextValue $ws.Range("E8") "  -1.60%  "
Set-TextValue $ws.Range("D9") "0.509"
Set-TextValue $ws.Range("E9") "  -2.36%  "
Set-TextValue $ws.Range("D10") "0.145"
Set-TextValue $ws.Range("E10") "  -2.64%  "
Set-TextValue $ws.Range("D11") "5.39"
Set-TextValue $ws.Range("E11") "  -1.10%  "
Set-TextValue $ws.Range("D12") "0.457"
Set-TextValue $ws.Range("E12") "  -3.28%  "
Set-TextValue $ws.Range("D13") "0.0000239"
Set-TextValue $ws.Range("E13") "  -3.53%  "
Set-TextValue $ws.Range("D14") "33.71"
Set-TextValue $ws.Range("E14") "  -3.49%  "
Set-TextValue $ws.Range("D15") "3.773.62"
Set-TextValue $ws.Range("E15") "  -1.35%  "
Set-TextValue $ws.Range("E16") "  +0.45%  "
Set-TextValue $ws.Range("D17") "3.245.42"
Set-TextValue $ws.Range("E17") "  -1.11%  "
Set-TextValue $ws.Range("D18") "62.856.92"
Set-TextValue $ws.Range("E18") "  -1.56%  "
Set-TextValue $ws.Range("D19") "6.73"
Set-TextValue $ws.Range("E19") "  -2.67%  "
Set-TextValue $ws.Range("D20") "468.89"
Set-TextValue $ws.Range("E20") "  -2.71%  "
Set-TextValue $ws.Range("D21") "13.81"
Set-TextValue $ws.Range("E21") "  -3.43%  "
Set-TextValue $ws.Range("D22") "0.719"
Set-TextValue $ws.Range("E22") "  -3.41%  "
Set-TextValue $ws.Range("E23") "  -3.89%  "
Set-TextValue $ws.Range("D24") "13.50"
Set-TextValue $ws.Range("E24") "  -0.35%  "
Set-TextValue $ws.Range("D25") "84.58"
Set-TextValue $ws.Range("E25") "  -0.12%  "
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.15%  "
Set-TextValue $ws.Range("D27") "2.71"
Set-TextValue $ws.Range("E27") "  -2.34%  "
Set-TextValue $ws.Range("E28") "  +0.05%  "
Set-TextValue $ws.Range("D29") "7.95"
Set-TextValue $ws.Range("E29") "  -4.08%  "
Set-TextValue $ws.Range("D30") "6.96"
Set-TextValue $ws.Range("E30") "  -4.42%  "
Set-TextValue $ws.Range("E31") "  -3.70%  "
Set-TextValue $ws.Range("D32") "27.91"
Set-TextValue $ws.Range("E32") "  -0.78%  "
Set-TextValue $ws.Range("E33") "  -5.56%  "
Set-TextValue $ws.Range("E34") "  -5.10%  "
Set-TextValue $ws.Range("E35") "  -3.70%  "
Set-TextValue $ws.Range("E36") "  -2.21%  "
Set-TextValue $ws.Range("D37") "51.82"
Set-TextValue $ws.Range("E37") "  -2.47%  "
Set-TextValue $ws.Range("D38") "0.0₃0716"
Set-TextValue $ws.Range("E38") "  -3.01%  "
Set-TextValue $ws.Range("D39") "0.0394"
Set-TextValue $ws.Range("E39") "  -1.09%  "
Set-TextValue $ws.Range("B40") "Maker"
Set-TextValue $ws.Range("C40") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D40") "3.037.52"
Set-TextValue $ws.Range("E40") "  +0.51%  "
Set-TextValue $ws.Range("B41") "Bittensor"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws.Range("D41") "420.84"
Set-TextValue $ws.Range("E41") "  -1.69%  "
Set-TextValue $ws.Range("D42") "0.118"
Set-TextValue $ws.Range("E42") "  +6.16%  "
Set-TextValue $ws.Range("D43") "8.13"
Set-TextValue $ws.Range("E43") "  -4.10%  "
Set-TextValue $ws.Range("E44") "  -5.93%  "
Set-TextValue $ws.Range("D45") "0.256"
Set-TextValue $ws.Range("E45") "  -5.20%  "
Set-TextValue $ws.Range("E46") "  -4.07%  "
Set-TextValue $ws.Range("B47") "USDe"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D47") "0.998"
Set-TextValue $ws.Range("E47") "  -0.09%  "
Set-TextValue $ws.Range("B48") "Arweave"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/7XWg41D1+arweave-ar"
Set-TextValue $ws.Range("D48") "35.62"
Set-TextValue $ws.Range("E48") "  +3.94%  "
Set-TextValue $ws.Range("D49") "126.18"
Set-TextValue $ws.Range("E49") "  +2.46%  "
Set-TextValue $ws.Range("D50") "25.76"
Set-TextValue $ws.Range("E50") "  -2.12%  "
Set-TextValue $ws.Range("D51") "0.112"
Set-TextValue $ws.Range("E51") "  -2.11%  "
